$wb = $excel.ActiveWorkbook

# --- REGCA1: remove the Iqmax / Iqmin parameter columns (U:V), shifting
#     the trailing ra / xs columns left into their place ---
$regca1 = $wb.Worksheets.Item("REGCA1")
$regca1.Range("U1:V2").Delete(-4159)

# --- REPCA1: insert a new "PLflag" parameter column right after "Fflag" ---
$repca1 = $wb.Worksheets.Item("REPCA1")
$repca1.Range("L1:L2").Insert(-4161)
$repca1.Range("L1").Value = "PLflag"
$repca1.Range("L2").Value = 0

$repca1.Activate()
